$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 234 (pushes existing rows 234-246 down to 235-247)
$ws.Rows.Item(234).Insert()

# Populate the newly inserted row 234 with the new record
$ws.Cells.Item(234, 1).Value = 3
$ws.Cells.Item(234, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(234, 3).Value = "Coquimbo"
$ws.Cells.Item(234, 4).Value = 44516
$ws.Cells.Item(234, 5).Value = 5
$ws.Cells.Item(234, 6).Value = 100112031
$ws.Cells.Item(234, 7).Value = "Poroto verde"
$ws.Cells.Item(234, 8).Value = "Magnum"
$ws.Cells.Item(234, 9).Value = "Primera"
$ws.Cells.Item(234, 10).Value = 73
$ws.Cells.Item(234, 11).Value = 40000
$ws.Cells.Item(234, 12).Value = 41000
$ws.Cells.Item(234, 13).Value = 40521
$ws.Cells.Item(234, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(234, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(234, 16).Value = 1621
$ws.Cells.Item(234, 17).Value = 25
$ws.Cells.Item(234, 18).Value = "Hortaliza"
